# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.829.82'
$ws.Range('E2').Value = '  +4.11%  '
$ws.Range('D3').Value = '2.630.07'
$ws.Range('E3').Value = '  +3.98%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.28'
$ws.Range('E6').Value = '  +2.21%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').Value = '2.628.88'
$ws.Range('E9').Value = '  +3.99%  '
$ws.Range('E10').Value = '  +12.54%  '
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.03'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('E14').Value = '  +4.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000186'
$ws.Range('E15').Value = '  +7.52%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '71.698.99'
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.54'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '2.649.79'
$ws.Range('E18').Value = '  +6.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '384.63'
$ws.Range('E19').Value = '  +5.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.96'
$ws.Range('E20').Value = '  +5.44%  '
$ws.Range('E21').Value = '  +4.57%  '
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.76'
$ws.Range('E23').Value = '  +3.05%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  +17.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.46'
$ws.Range('E25').Value = '  +6.32%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('E27').Value = '  +10.12%  '
$ws.Range('D28').Value = '2.767.23'
$ws.Range('E28').Value = '  +4.33%  '
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').Value = '0.0₃0964'
$ws.Range('E30').Value = '  +8.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '544.60'
$ws.Range('E31').Value = '  +5.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.05'
$ws.Range('E32').Value = '  +3.38%  '
$ws.Range('E33').Value = '  +7.76%  '
$ws.Range('E34').Value = '  +2.99%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.01'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('E40').Value = '  +6.17%  '
$ws.Range('E41').Value = '  +7.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.63'
$ws.Range('E42').Value = '  +10.78%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.03'
$ws.Range('E44').Value = '  +4.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.23'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '151.05'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.535'
$ws.Range('E49').Value = '  +3.72%  '
$ws.Range('E50').Value = '  +7.28%  '
$ws.Range('D51').Value = '0.0₆0264'
$ws.Range('E51').Value = '  +5.19%  '
